$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capture the note text used by the previous row so the new row reuses
# the same shared string (same as the existing pattern in the sheet).
$prevNote = $ws.Range("B60").Value2

# Copy the formatting (borders/fill/alignment) of row 60 down to the new
# row 61, matching the style already used for the data rows.
$ws.Range("A60:B60").Copy() | Out-Null
$ws.Range("A61:B61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Add the new day's date entry.
$ws.Range("A61").Value = "16-11-2025"
$ws.Range("B61").Value = $prevNote
